# Generate Report for Handoff
# Regenerate the localization-status report cells: new handoff id (UUID),
# new content hashes, and refreshed timestamps. Mirrors what the CI
# report-generation job produces on each run.

$wb = $excel.ActiveWorkbook

$oldId  = "fc066dd0-e8e5-4755-a690-7d0357d4893d"
$newId  = "80aa41f0-5531-48e8-8ca9-9e525a3dee92"
$oldZh  = "4c63327ace79fd7f3a03ca82e028b634c0473465"
$newZh  = "85403578801eaa6a87a245f336e64380eb140db0"
$oldDe  = "4c63327ace79fd7f3a03ca82e028b634c0473465"
$newDe  = "85403578801eaa6a87a245f336e64380eb140db0"

$newMdName    = "$newId.md"
$newZhXlfName = "$newId.$newZh.zh-cn.xlf"
$newDeXlfName = "$newId.$newDe.de-de.xlf"

# NOTE: the hyperlink *targets* (the github blob URLs, pinned to an old
# commit SHA and the old file name) are left exactly as they were -- only
# the visible label (TextToDisplay) and the backing shared-string text are
# refreshed to the new handoff id/hash. That matches the source report
# generator, which rewrites cell text without touching the relationship
# targets already recorded for this row.
$mdAddr    = "https://github.com/OpenLocalizationTest/oltest/blob/9df54e1d6616c7ae3124776ed06647a18a3f5aec/e2e/$oldId.md"
$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/541a3634bc12df1f0d9706390289a99fda426f4d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldZh.zh-cn.xlf"
$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f0bf269c52feb28b7a58e630aa8f403e764c8f5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldDe.de-de.xlf"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddr, "", "", $newMdName)
$ws.Range("D2").Value2 = "2016-47-20 08:47:48"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddr, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdAddr, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $zhXlfAddr, "", "", $newZhXlfName)
$ws.Range("E2").Value2 = "2016-03-20 08:47:45"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddr, "", "", $newMdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdAddr, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $deXlfAddr, "", "", $newDeXlfName)
$ws.Range("E2").Value2 = "2016-03-20 08:47:48"
